# ----------------------------------------------------------------------------
# Scheduled market-data refresh for Asura_Profits workbook (Universalis price pull).
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ,
# LevePriceNQ / LevePriceHQ and the derived LeveProfitNQ / LeveProfitHQ columns
# (H:N) for the leves whose underlying item prices moved since the last run.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ==================== ALC ====================
$ws = $wb.Worksheets.Item("ALC")
# Row 64: "Forged from the Void" (Void Glue)
$ws.Range("I64").Value = 3128.3572
$ws.Range("J64").Value = 4255.5
$ws.Range("K64").Value = 3128.3572
$ws.Range("L64").Value = 4255.5
$ws.Range("M64").Value = -2880.3572
$ws.Range("N64").Value = -4751.5
# Row 67: "Dodging the Draft (L)" (Void Glue)
$ws.Range("I67").Value = 3128.3572
$ws.Range("J67").Value = 4255.5
$ws.Range("K67").Value = 3128.3572
$ws.Range("L67").Value = 4255.5
$ws.Range("M67").Value = -2270.3572
$ws.Range("N67").Value = -5971.5
# Row 76: "Warding Off Temptation" (Enchanted Hardsilver Ink)
$ws.Range("H76").Value = 4637.5
$ws.Range("I76").Value = 5133.3335
$ws.Range("K76").Value = 5133.3335
$ws.Range("M76").Value = -4818.3335
# Row 79: "The Garden of Arcane Delights (L)" (Enchanted Hardsilver Ink)
$ws.Range("H79").Value = 4637.5
$ws.Range("I79").Value = 5133.3335
$ws.Range("K79").Value = 5133.3335
$ws.Range("M79").Value = -4041.3335

# ==================== ARM ====================
$ws = $wb.Worksheets.Item("ARM")
# Row 63: "Rivets Run through It" (Mythrite Rivets)
$ws.Range("H63").Value = 3848.7693
$ws.Range("I63").Value = 3336.1667
$ws.Range("K63").Value = 3336.1667
$ws.Range("M63").Value = -2650.1667
# Row 66: "A Riveting Revival (L)" (Mythrite Rivets)
$ws.Range("H66").Value = 3848.7693
$ws.Range("I66").Value = 3336.1667
$ws.Range("K66").Value = 16680.8335
$ws.Range("M66").Value = -13248.8335
# Row 132: "Don't Bore Me, Ore Me" (Mountain Chromite Ingot)
$ws.Range("H132").Value = 1112804.9
$ws.Range("I132").Value = 1539268.8
$ws.Range("K132").Value = 4617806.4
$ws.Range("M132").Value = -4615276.4

# ==================== BSM ====================
$ws = $wb.Worksheets.Item("BSM")
# Row 105: "Ingot to Wing It" (Molybdenum Ingot)
$ws.Range("H105").Value = 2364.1
$ws.Range("I105").Value = 2382.3333
$ws.Range("J105").Value = 2200
$ws.Range("K105").Value = 2382.3333
$ws.Range("L105").Value = 2200
$ws.Range("M105").Value = -635.3332999999998
$ws.Range("N105").Value = -5694
# Row 107: "The Gold Experience" (Deepgold Nugget)
$ws.Range("H107").Value = 48996.453
$ws.Range("I107").Value = 53565.1
$ws.Range("J107").Value = 3310
$ws.Range("K107").Value = 53565.1
$ws.Range("L107").Value = 3310
$ws.Range("M107").Value = -51645.1
$ws.Range("N107").Value = -7150
# Row 134: "Ruthenium Supremium" (Ruthenium Ingot)
$ws.Range("H134").Value = 515630.28
$ws.Range("I134").Value = 802295.6
$ws.Range("K134").Value = 2406886.8
$ws.Range("M134").Value = -2404351.8

# ==================== CRP ====================
$ws = $wb.Worksheets.Item("CRP")
# Row 8: "Bows for the Boys" (Maple Longbow)
$ws.Range("H8").Value = 3138.3333
$ws.Range("I8").Value = 1500
$ws.Range("J8").Value = 3957.5
$ws.Range("K8").Value = 1500
$ws.Range("L8").Value = 3957.5
$ws.Range("M8").Value = -1360
$ws.Range("N8").Value = -4237.5
# Row 25: "Bowing to Necessity" (Ash Shortbow)
$ws.Range("H25").Value = 1000000000
$ws.Range("I25").Value = 1000000000
$ws.Range("K25").Value = 1000000000
$ws.Range("M25").Value = -999999826
# Row 62: "Splinter in the Sewers" (Cedar Lumber)
$ws.Range("H62").Value = 60623.332
$ws.Range("I62").Value = 85868.336
$ws.Range("J62").Value = 10133.333
$ws.Range("K62").Value = 85868.336
$ws.Range("L62").Value = 10133.333
$ws.Range("M62").Value = -85244.336
$ws.Range("N62").Value = -11381.333
# Row 65: "The Lumber of Their Discontent (L)" (Cedar Lumber)
$ws.Range("H65").Value = 60623.332
$ws.Range("I65").Value = 85868.336
$ws.Range("J65").Value = 10133.333
$ws.Range("K65").Value = 429341.68
$ws.Range("L65").Value = 50666.665
$ws.Range("M65").Value = -426221.68
$ws.Range("N65").Value = -56906.665
# Row 86: "Birch, Please" (Birch Lumber)
$ws.Range("N86").ClearContents()
$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1877
# Row 89: "Built This City on Blocks and Soul (L)" (Birch Lumber)
$ws.Range("N89").ClearContents()
$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -9384
# Row 99: "O Pine" (Pine Lumber)
$ws.Range("H99").Value = 2554.3704
$ws.Range("I99").Value = 2555.7144
$ws.Range("J99").Value = 2549.6667
$ws.Range("K99").Value = 2555.7144
$ws.Range("L99").Value = 2549.6667
$ws.Range("M99").Value = -1057.7144
$ws.Range("N99").Value = -5545.6667
# Row 126: "A Better Conductor" (Red Pine Lumber)
$ws.Range("H126").Value = 2554.3704
$ws.Range("I126").Value = 2555.7144
$ws.Range("J126").Value = 2549.6667
$ws.Range("K126").Value = 7667.1432
$ws.Range("L126").Value = 7649.000100000001
$ws.Range("M126").Value = -5197.1432
$ws.Range("N126").Value = -12589.0001
# Row 134: "Wood You Be Quiet" (Ceiba Lumber)
$ws.Range("H134").Value = 1653.4762
$ws.Range("I134").Value = 1510.6
$ws.Range("K134").Value = 4531.799999999999
$ws.Range("M134").Value = -1996.799999999999

# ==================== GSM ====================
$ws = $wb.Worksheets.Item("GSM")
# Row 70: "Sky Is the Limit" (Mythrite Ingot)
$ws.Range("H70").Value = 8754.25
$ws.Range("I70").Value = 8504
$ws.Range("K70").Value = 8504
$ws.Range("M70").Value = -8234
# Row 73: "Hulls of Broken Dreams (L)" (Mythrite Ingot)
$ws.Range("H73").Value = 8754.25
$ws.Range("I73").Value = 8504
$ws.Range("K73").Value = 8504
$ws.Range("M73").Value = -7568
# Row 80: "Needs More Prayerbell" (Hardsilver Ingot)
$ws.Range("H80").Value = 2573.2666
$ws.Range("I80").Value = 2355.5557
$ws.Range("J80").Value = 2899.8333
$ws.Range("K80").Value = 2355.5557
$ws.Range("L80").Value = 2899.8333
$ws.Range("M80").Value = -1357.5557
$ws.Range("N80").Value = -4895.8333
# Row 83: "With a Noise That Reaches Heaven (L)" (Hardsilver Ingot)
$ws.Range("H83").Value = 2573.2666
$ws.Range("I83").Value = 2355.5557
$ws.Range("J83").Value = 2899.8333
$ws.Range("K83").Value = 11777.7785
$ws.Range("L83").Value = 14499.1665
$ws.Range("M83").Value = -6785.7785
$ws.Range("N83").Value = -24483.1665
# Row 102: "Put the Metal to the Peddle" (Durium Ingot)
$ws.Range("H102").Value = 2128.7812
$ws.Range("I102").Value = 2124.652
$ws.Range("J102").Value = 2139.3333
$ws.Range("K102").Value = 2124.652
$ws.Range("L102").Value = 2139.3333
$ws.Range("M102").Value = -502.652
$ws.Range("N102").Value = -5383.3333
# Row 126: "Gold Rush Order" (Phrygian Gold Ingot)
$ws.Range("H126").Value = 2914.8572
$ws.Range("I126").Value = 2617.7778
$ws.Range("J126").Value = 3449.6
$ws.Range("K126").Value = 7853.3334
$ws.Range("L126").Value = 10348.8
$ws.Range("M126").Value = -5383.3334
$ws.Range("N126").Value = -15288.8

# ==================== LTW ====================
$ws = $wb.Worksheets.Item("LTW")
# Row 24: "On Their Feet Again" (Hard Leather Espadrilles)
$ws.Range("M24").ClearContents()
$ws.Range("H24").Value = 7995
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 7995
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 7995
$ws.Range("N24").Value = -8681

# ==================== WVR ====================
$ws = $wb.Worksheets.Item("WVR")
# Row 20: "Read the Fine Print" (Cotton Shepherd's Tunic)
$ws.Range("H20").Value = 22507.875
$ws.Range("J20").Value = 13343.833
$ws.Range("L20").Value = 13343.833
$ws.Range("N20").Value = -13823.833
# Row 122: "Heavy Armoire" (Dark Hempen Cloth)
$ws.Range("H122").Value = 31254970
$ws.Range("I122").Value = 41669030
$ws.Range("J122").Value = 12777.5
$ws.Range("K122").Value = 125007090
$ws.Range("L122").Value = 38332.5
$ws.Range("M122").Value = -125004640
$ws.Range("N122").Value = -43232.5
# Row 132: "Comfy Cabins" (Snow Cotton Cloth)
$ws.Range("H132").Value = 2319.0588
$ws.Range("I132").Value = 1759.45
$ws.Range("K132").Value = 5278.35
$ws.Range("M132").Value = -2748.35
